$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 10:46"

# Row 8: Alemania
$ws.Range("B8").Value = 30081
$ws.Range("C8").Value = 1025
$ws.Range("D8").Value = 453
$ws.Range("E8").Value = 29498
$ws.Range("F8").Value = 23
$ws.Range("G8").Value = 7
$ws.Range("H8").Value = 130

# Row 12: Suiza
$ws.Range("B12").Value = 8795
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 131
$ws.Range("E12").Value = 8547
$ws.Range("F12").Value = 141
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 117

# Row 15: Austria
$ws.Range("B15").Value = 4668
$ws.Range("C15").Value = 194
$ws.Range("D15").Value = 9
$ws.Range("E15").Value = 4634
$ws.Range("F15").Value = 16
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = 25

# Row 34: Polonia
$ws.Range("B34").Value = 774
$ws.Range("C34").Value = 25
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 764
$ws.Range("F34").Value = 3
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 9

# Row 89: Sri Lanka
$ws.Range("A89").Value = "Sri Lanka"
$ws.Range("B89").Value = 100
$ws.Range("C89").Value = 3
$ws.Range("D89").Value = 2
$ws.Range("E89").Value = 98
$ws.Range("F89").Value = 2
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 0

# Row 90: Burkina Faso
$ws.Range("A90").Value = "Burkina Faso"
$ws.Range("B90").Value = 99
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 5
$ws.Range("E90").Value = 90
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 4

# Row 102: Georgia
$ws.Range("A102").Value = "Georgia"
$ws.Range("B102").Value = 67
$ws.Range("C102").Value = 6
$ws.Range("D102").Value = 9
$ws.Range("E102").Value = 58
$ws.Range("F102").Value = 1
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 0

# Row 103: Camerun
$ws.Range("A103").Value = "Camerun"
$ws.Range("B103").Value = 66
$ws.Range("C103").Value = 10
$ws.Range("D103").Value = 2
$ws.Range("E103").Value = 64
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 0

# Row 105: Estado de Palestina
$ws.Range("B105").Value = 60
$ws.Range("C105").Value = 1
$ws.Range("D105").Value = 17
$ws.Range("E105").Value = 43
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 0

# Row 115: Banglades
$ws.Range("A115").Value = "Banglades"
$ws.Range("B115").Value = 39
$ws.Range("C115").Value = 6
$ws.Range("D115").Value = 5
$ws.Range("E115").Value = 30
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 4

# Row 116: Ruanda
$ws.Range("A116").Value = "Ruanda"
$ws.Range("B116").Value = 36
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 0
$ws.Range("E116").Value = 36
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 0

# Row 117: Mayotte
$ws.Range("A117").Value = "Mayotte"
$ws.Range("B117").Value = 36
$ws.Range("C117").Value = 12
$ws.Range("D117").Value = 0
$ws.Range("E117").Value = 36
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 0

# Row 118: Mauricio
$ws.Range("A118").Value = "Mauricio"
$ws.Range("B118").Value = 36
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 0
$ws.Range("E118").Value = 34
$ws.Range("F118").Value = 1
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 2

# Row 135: Islas Virgenes de los Estados Unidos
$ws.Range("A135").Value = "Islas Virgenes de los Estados Unidos"

# Row 136: Madagascar
$ws.Range("A136").Value = "Madagascar"
$ws.Range("B136").Value = 17
$ws.Range("C136").Value = 5
$ws.Range("D136").Value = 0
$ws.Range("E136").Value = 17
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0

# Row 137: Barbados
$ws.Range("A137").Value = "Barbados"
$ws.Range("B137").Value = 17
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 0
$ws.Range("E137").Value = 17
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 0

# Row 142: Etiopia
$ws.Range("A142").Value = "Etiopia"
$ws.Range("B142").Value = 12
$ws.Range("C142").Value = 1
$ws.Range("D142").Value = 0
$ws.Range("E142").Value = 12
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 0

# Row 143: Tanzania
$ws.Range("A143").Value = "Tanzania"
$ws.Range("B143").Value = 12
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 0
$ws.Range("E143").Value = 12
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 0

# Row 145: Nueva Caledonia
$ws.Range("A145").Value = "Nueva Caledonia"
$ws.Range("B145").Value = 10
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 0
$ws.Range("E145").Value = 10
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 0

# Row 146: Mongolia
$ws.Range("A146").Value = "Mongolia"
$ws.Range("B146").Value = 10
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 0
$ws.Range("E146").Value = 10
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 0

# Row 159: Groenlandia
$ws.Range("A159").Value = "Groenlandia"

# Row 160: Congo
$ws.Range("A160").Value = "Congo"

# Row 161: Suazilandia
$ws.Range("A161").Value = "Suazilandia"

# Row 162: Guinea
$ws.Range("A162").Value = "Guinea"

# Row 164: Bahamas
$ws.Range("A164").Value = "Bahamas"

# Row 166: Republica del Chad
$ws.Range("A166").Value = "Republica del Chad"
$ws.Range("B166").Value = 3
$ws.Range("C166").Value = 1
$ws.Range("D166").Value = 0
$ws.Range("E166").Value = 3
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0

# Row 171: Cabo Verde
$ws.Range("A171").Value = "Cabo Verde"

# Row 172: Santa Lucia
$ws.Range("A172").Value = "Santa Lucia"

# Row 173: Zambia
$ws.Range("A173").Value = "Zambia"

# Row 174: Republica de Africa Central
$ws.Range("A174").Value = "Republica de Africa Central"

# Row 175: Republica de Yibuti
$ws.Range("A175").Value = "Republica de Yibuti"

# Row 176: San Bartolome
$ws.Range("A176").Value = "San Bartolome"
$ws.Range("B176").Value = 3
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 3
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 0

# Row 177: Zimbabue
$ws.Range("A177").Value = "Zimbabue"
$ws.Range("B177").Value = 3
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 0
$ws.Range("E177").Value = 2
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 1

# Row 178: Birmania
$ws.Range("A178").Value = "Birmania"

# Row 179: Nicaragua
$ws.Range("A179").Value = "Nicaragua"

# Row 181: Butan
$ws.Range("A181").Value = "Butan"

# Row 182: Dominica
$ws.Range("A182").Value = "Dominica"

# Row 183: Mauritania
$ws.Range("A183").Value = "Mauritania"

# Row 184: San Martin (Parte Holandesa)
$ws.Range("A184").Value = "San Martin (Parte Holandesa)"
